# cryptos.xlsx refresh -- Tue Jan 10 23:46:51 UTC 2023 symbol-list update
# Re-pull of Price (D) and Volume(1h) (E) for the tracked coins.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $rng = $ws.Range($addr)
    # Force text entry so numeric-looking strings (prices, percentages)
    # stay literal text instead of being parsed into numbers/dates,
    # then drop back to the Normal style so no formatting lingers.
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "277.06"
Set-TextCell "E2" "1.47%"
Set-TextCell "D3" "27.16"
Set-TextCell "E3" "1.70%"
Set-TextCell "D4" "4.888"
Set-TextCell "E4" "-0.34%"
Set-TextCell "D5" "0.06421"
Set-TextCell "E5" "1.56%"
Set-TextCell "D6" "6.946"
Set-TextCell "E6" "0.50%"
Set-TextCell "E7" "-8.57%"
Set-TextCell "E8" "-0.44%"
Set-TextCell "D9" "0.1519"
Set-TextCell "E9" "3.95%"
Set-TextCell "D10" "0.05025"
Set-TextCell "E10" "-2.15%"
Set-TextCell "D11" "0.07505"
Set-TextCell "E11" "1.56%"
Set-TextCell "D12" "0.02919"
Set-TextCell "E12" "-6.67%"
Set-TextCell "D13" "0.09000"
Set-TextCell "E13" "-0.47%"
Set-TextCell "D14" "0.001573"
Set-TextCell "E14" "0.77%"
Set-TextCell "D15" "0.0006404"
Set-TextCell "E15" "1.42%"
Set-TextCell "D16" "0.005756"
Set-TextCell "E16" "-4.58%"
Set-TextCell "E17" "0.08%"
Set-TextCell "D18" "3.316"
Set-TextCell "E18" "-1.17%"
Set-TextCell "E19" "-0.47%"
Set-TextCell "E20" "-0.94%"
Set-TextCell "E21" "1.60%"
Set-TextCell "E22" "0.36%"
Set-TextCell "D23" "0.04432"
Set-TextCell "E23" "1.43%"
Set-TextCell "D24" "0.001174"
Set-TextCell "E24" "-0.17%"
Set-TextCell "E25" "5.00%"
Set-TextCell "E27" "13.98%"
Set-TextCell "D40" "0.04146"
Set-TextCell "E40" "2.86%"
Set-TextCell "D41" "0.006815"
Set-TextCell "E41" "2.84%"
Set-TextCell "D42" "0.1177"
Set-TextCell "E42" "0.91%"
Set-TextCell "E43" "13.93%"
Set-TextCell "D44" "0.01173"
Set-TextCell "E44" "-4.20%"
Set-TextCell "D45" "0.00005210"
Set-TextCell "E45" "-2.03%"
Set-TextCell "E46" "-36.92%"
Set-TextCell "E47" "-22.19%"
